# Apply PO export changes: new supplier info, new PO details, and new line items.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("For End-User")

# --- Header block: supplier / PO no / address / date ---
$ws1.Range("C7").Value = "Alprops Management  Inc. ( Montevista Hot Spring and Conference Resort)"
$ws1.Range("F7").Value = "2022-03-0005"
$ws1.Range("C8").Value = "Barangay Pansol Calamba Laguna"
$ws1.Range("F8").Value = "April 08, 2022"

# Mode of Procurement value cleared
$ws1.Range("E10").Value = ""

# --- Line item rows (17-19) ---
# Row 17
$ws1.Range("A17").Value = "S3607"
$ws1.Range("B17").Value = "pack"
$ws1.Range("C17").Value = "Specialty Paper, A4`nTEST7"
$ws1.Range("E17").Value = 1
$ws1.Range("F17").Value = 60
$ws1.Range("G17").Value = 60
$ws1.Rows.Item(17).RowHeight = 30

# Row 18
$ws1.Range("A18").Value = "S3634"
$ws1.Range("B18").Value = "piece"
$ws1.Range("C18").Value = "Flasher`nTEST7"
$ws1.Range("E18").Value = 2
$ws1.Range("F18").Value = 400
$ws1.Range("G18").Value = 800
$ws1.Rows.Item(18).RowHeight = 30

# Row 19
$ws1.Range("A19").Value = "S3642"
$ws1.Range("B19").Value = "lot"
$ws1.Range("C19").Value = "Parking Fee`nTEST7"
$ws1.Range("E19").Value = 3
$ws1.Range("F19").Value = 150
$ws1.Range("G19").Value = 450
$ws1.Rows.Item(19).RowHeight = 30

# --- Totals block ---
$ws1.Range("G48").Value = 0
$ws1.Range("A48").Value = "(Total Amount in Words)    pesos only"
